# death_causes.xlsx — "modified all data to facilitate df import"
#
# 1. Three missing years (1944-1946) are inserted into the time series
#    (between 1943 and 1947), pushing all subsequent data rows down by 3
#    (old row 47..116 -> new row 50..119). Only the Year column is known
#    for these new rows, so the rest of the row stays blank.
# 2. The "Death rates" column header (column B) is renamed to "deaths".
# 3. The emc120000_ defined name is updated to the new used range.
# 4. The view is scrolled down to where the data now is, with a new
#    active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 3 rows for the missing years 1944, 1945, 1946 --------------
# Old row 47 (year 1947) is the first row after 1943 (row 46); push it (and
# everything below it) down by 3 rows so the new rows can hold 1944-1946.
$ws.Rows("47:49").Insert()

$ws.Range("A47").Value = 1944
$ws.Range("A48").Value = 1945
$ws.Range("A49").Value = 1946

# --- 2. Rename column B header from "Death rates" to "deaths" -------------
$ws.Range("B1").Value = "deaths"

# --- 3. Update the emc120000_ defined name to Sheet1!$A$1:$J$119 ----------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!emc120000_") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$J`$119"
    }
}

# --- 4. Update the sheet view: scroll to row 83 and select N110 -----------
$excel.ActiveWindow.ScrollRow = 83
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N110").Select()
